$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.0407475
$ws.Range("H2").Value = 4.081495
$ws.Range("I2").Value = 0.007689531190315315
$ws.Range("J2").Value = 0.005145434208836267
$ws.Range("M2").Value = 181.556244
$ws.Range("N2").Value = 544.668732
$ws.Range("O2").Value = 0.393453292404907
$ws.Range("P2").Value = 0.3935455037432071
$ws.Range("Q2").Value = 370.51045105239
$ws.Range("R2").Value = 2223.06270631434
$ws.Range("S2").Value = 0.003025471363879784
$ws.Range("T2").Value = 0.002024962497693999
# Row 3
$ws.Range("G3").Value = 2.0407475
$ws.Range("H3").Value = 4.081495
$ws.Range("I3").Value = 0.007689531190315315
$ws.Range("J3").Value = 0.005145434208836267
$ws.Range("M3").Value = 0.324361
$ws.Range("N3").Value = 0.648722
$ws.Range("O3").Value = 0.0007029276469155644
$ws.Range("P3").Value = 0.0004687282586276696
$ws.Range("Q3").Value = 0.6619388998475001
$ws.Range("R3").Value = 2.64775559939
$ws.Range("S3").Value = 0.000005405184065492183
$ws.Range("T3").Value = 0.000002411810416591064
# Row 4
$ws.Range("G4").Value = 2.0407475
$ws.Range("H4").Value = 4.081495
$ws.Range("I4").Value = 0.007689531190315315
$ws.Range("J4").Value = 0.005145434208836267
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 155.929759
$ws.Range("N4").Value = 467.789277
$ws.Range("O4").Value = 0.3379177477501335
$ws.Range("P4").Value = 0.3379969435488647
$ws.Range("Q4").Value = 318.2132658548526
$ws.Range("R4").Value = 1909.279595129115
$ws.Range("S4").Value = 0.002598429061085754
$ws.Range("T4").Value = 0.001739141035818429
# Row 5
$ws.Range("G5").Value = 2.0407475
$ws.Range("H5").Value = 4.081495
$ws.Range("I5").Value = 0.007689531190315315
$ws.Range("J5").Value = 0.005145434208836267
$ws.Range("M5").Value = 123.632576
$ws.Range("N5").Value = 370.897728
$ws.Range("O5").Value = 0.2679260321980438
$ws.Range("P5").Value = 0.2679888244493004
$ws.Range("Q5").Value = 252.3028703905601
$ws.Range("R5").Value = 1513.81722234336
$ws.Range("S5").Value = 0.002060225581284283
$ws.Range("T5").Value = 0.001378918864907247
# Row 6
$ws.Range("I6").Value = 0.4072821437310581
$ws.Range("J6").Value = 0.4087980313366845
$ws.Range("M6").Value = 181.556244
$ws.Range("N6").Value = 544.668732
$ws.Range("O6").Value = 0.393453292404907
$ws.Range("P6").Value = 0.3935455037432071
$ws.Range("Q6").Value = 19624.3811286486
$ws.Range("R6").Value = 176619.4301578374
$ws.Range("S6").Value = 0.1602465003887134
$ws.Range("T6").Value = 0.1608806271716269
# Row 7
$ws.Range("I7").Value = 0.4072821437310581
$ws.Range("J7").Value = 0.4087980313366845
$ws.Range("M7").Value = 0.324361
$ws.Range("N7").Value = 0.648722
$ws.Range("O7").Value = 0.0007029276469155644
$ws.Range("P7").Value = 0.0004687282586276696
$ws.Range("Q7").Value = 35.06012102381667
$ws.Range("R7").Value = 210.3607261429
$ws.Range("S7").Value = 0.0002862898789235993
$ws.Range("T7").Value = 0.0001916151893588636
# Row 8
$ws.Range("I8").Value = 0.4072821437310581
$ws.Range("J8").Value = 0.4087980313366845
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 155.929759
$ws.Range("N8").Value = 467.789277
$ws.Range("O8").Value = 0.3379177477501335
$ws.Range("P8").Value = 0.3379969435488647
$ws.Range("Q8").Value = 16854.41906318752
$ws.Range("R8").Value = 151689.7715686877
$ws.Range("S8").Value = 0.1376278647084453
$ws.Range("T8").Value = 0.1381724851205924
# Row 9
$ws.Range("I9").Value = 0.4072821437310581
$ws.Range("J9").Value = 0.4087980313366845
$ws.Range("M9").Value = 123.632576
$ws.Range("N9").Value = 370.897728
$ws.Range("O9").Value = 0.2679260321980438
$ws.Range("P9").Value = 0.2679888244493004
$ws.Range("Q9").Value = 13363.42247386773
$ws.Range("R9").Value = 120270.8022648096
$ws.Range("S9").Value = 0.1091214887549758
$ws.Range("T9").Value = 0.1095533038551064
# Row 10
$ws.Range("G10").Value = 49.59263000000001
$ws.Range("H10").Value = 148.77789
$ws.Range("I10").Value = 0.18686489886415
$ws.Range("J10").Value = 0.1875604024320694
$ws.Range("M10").Value = 181.556244
$ws.Range("N10").Value = 544.668732
$ws.Range("O10").Value = 0.393453292404907
$ws.Range("P10").Value = 0.3935455037432071
$ws.Range("Q10").Value = 9003.85163288172
$ws.Range("R10").Value = 81034.66469593548
$ws.Range("S10").Value = 0.07352260969300979
$ws.Range("T10").Value = 0.07381355305740742
# Row 11
$ws.Range("G11").Value = 49.59263000000001
$ws.Range("H11").Value = 148.77789
$ws.Range("I11").Value = 0.18686489886415
$ws.Range("J11").Value = 0.1875604024320694
$ws.Range("M11").Value = 0.324361
$ws.Range("N11").Value = 0.648722
$ws.Range("O11").Value = 0.0007029276469155644
$ws.Range("P11").Value = 0.0004687282586276696
$ws.Range("Q11").Value = 16.08591505943
$ws.Range("R11").Value = 96.51549035658002
$ws.Range("S11").Value = 0.0001313525036496919
$ws.Range("T11").Value = 0.00008791486081948882
# Row 12
$ws.Range("G12").Value = 49.59263000000001
$ws.Range("H12").Value = 148.77789
$ws.Range("I12").Value = 0.18686489886415
$ws.Range("J12").Value = 0.1875604024320694
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 155.929759
$ws.Range("N12").Value = 467.789277
$ws.Range("O12").Value = 0.3379177477501335
$ws.Range("P12").Value = 0.3379969435488647
$ws.Range("Q12").Value = 7732.966844076172
$ws.Range("R12").Value = 69596.70159668554
$ws.Range("S12").Value = 0.06314496575773006
$ws.Range("T12").Value = 0.06339484275283451
# Row 13
$ws.Range("G13").Value = 49.59263000000001
$ws.Range("H13").Value = 148.77789
$ws.Range("I13").Value = 0.18686489886415
$ws.Range("J13").Value = 0.1875604024320694
$ws.Range("M13").Value = 123.632576
$ws.Range("N13").Value = 370.897728
$ws.Range("O13").Value = 0.2679260321980438
$ws.Range("P13").Value = 0.2679888244493004
$ws.Range("Q13").Value = 6131.264597514882
$ws.Range("R13").Value = 55181.38137763393
$ws.Range("S13").Value = 0.05006597090976046
$ws.Range("T13").Value = 0.05026409176100799
# Row 14
$ws.Range("G14").Value = 0.9116095
$ws.Range("H14").Value = 1.823219
$ws.Range("I14").Value = 0.003434942188407801
$ws.Range("J14").Value = 0.002298484602529281
$ws.Range("M14").Value = 181.556244
$ws.Range("N14").Value = 544.668732
$ws.Range("O14").Value = 0.393453292404907
$ws.Range("P14").Value = 0.3935455037432071
$ws.Range("Q14").Value = 165.508396814718
$ws.Range("R14").Value = 993.0503808883079
$ws.Range("S14").Value = 0.001351489313249566
$ws.Range("T14").Value = 0.0009045582807483912
# Row 15
$ws.Range("G15").Value = 0.9116095
$ws.Range("H15").Value = 1.823219
$ws.Range("I15").Value = 0.003434942188407801
$ws.Range("J15").Value = 0.002298484602529281
$ws.Range("M15").Value = 0.324361
$ws.Range("N15").Value = 0.648722
$ws.Range("O15").Value = 0.0007029276469155644
$ws.Range("P15").Value = 0.0004687282586276696
$ws.Range("Q15").Value = 0.2956905690295
$ws.Range("R15").Value = 1.182762276118
$ws.Range("S15").Value = 0.000002414515829788495
$ws.Range("T15").Value = 0.000001077364685226061
# Row 16
$ws.Range("G16").Value = 0.9116095
$ws.Range("H16").Value = 1.823219
$ws.Range("I16").Value = 0.003434942188407801
$ws.Range("J16").Value = 0.002298484602529281
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 155.929759
$ws.Range("N16").Value = 467.789277
$ws.Range("O16").Value = 0.3379177477501335
$ws.Range("P16").Value = 0.3379969435488647
$ws.Range("Q16").Value = 142.1470496371105
$ws.Range("R16").Value = 852.882297822663
$ws.Range("S16").Value = 0.001160727927958679
$ws.Range("T16").Value = 0.0007768807704490241
# Row 17
$ws.Range("G17").Value = 0.9116095
$ws.Range("H17").Value = 1.823219
$ws.Range("I17").Value = 0.003434942188407801
$ws.Range("J17").Value = 0.002298484602529281
$ws.Range("M17").Value = 123.632576
$ws.Range("N17").Value = 370.897728
$ws.Range("O17").Value = 0.2679260321980438
$ws.Range("P17").Value = 0.2679888244493004
$ws.Range("Q17").Value = 112.704630791072
$ws.Range("R17").Value = 676.227784746432
$ws.Range("S17").Value = 0.0009203104313697674
$ws.Range("T17").Value = 0.0006159681866466395
# Row 18
$ws.Range("G18").Value = 97.62255466666666
$ws.Range("H18").Value = 292.867664
$ws.Range("I18").Value = 0.3678415281594588
$ws.Range("J18").Value = 0.3692106193949926
$ws.Range("M18").Value = 181.556244
$ws.Range("N18").Value = 544.668732
$ws.Range("O18").Value = 0.393453292404907
$ws.Range("P18").Value = 0.3935455037432071
$ws.Range("Q18").Value = 17723.98435496467
$ws.Range("R18").Value = 159515.859194682
$ws.Range("S18").Value = 0.1447284603375914
$ws.Range("T18").Value = 0.1453011791971439
# Row 19
$ws.Range("G19").Value = 97.62255466666666
$ws.Range("H19").Value = 292.867664
$ws.Range("I19").Value = 0.3678415281594588
$ws.Range("J19").Value = 0.3692106193949926
$ws.Range("M19").Value = 0.324361
$ws.Range("N19").Value = 0.648722
$ws.Range("O19").Value = 0.0007029276469155644
$ws.Range("P19").Value = 0.0004687282586276696
$ws.Range("Q19").Value = 31.66494945423467
$ws.Range("R19").Value = 189.989696725408
$ws.Range("S19").Value = 0.0002585659798269537
$ws.Range("T19").Value = 0.0001730594506958582
# Row 20
$ws.Range("G20").Value = 97.62255466666666
$ws.Range("H20").Value = 292.867664
$ws.Range("I20").Value = 0.3678415281594588
$ws.Range("J20").Value = 0.3692106193949926
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 155.929759
$ws.Range("N20").Value = 467.789277
$ws.Range("O20").Value = 0.3379177477501335
$ws.Range("P20").Value = 0.3379969435488647
$ws.Range("Q20").Value = 15222.26142213766
$ws.Range("R20").Value = 137000.3527992389
$ws.Range("S20").Value = 0.1243001807246116
$ws.Range("T20").Value = 0.1247920608812907
# Row 21
$ws.Range("G21").Value = 97.62255466666666
$ws.Range("H21").Value = 292.867664
$ws.Range("I21").Value = 0.3678415281594588
$ws.Range("J21").Value = 0.3692106193949926
$ws.Range("M21").Value = 123.632576
$ws.Range("N21").Value = 370.897728
$ws.Range("O21").Value = 0.2679260321980438
$ws.Range("P21").Value = 0.2679888244493004
$ws.Range("Q21").Value = 12069.32790914082
$ws.Range("R21").Value = 108623.9511822674
$ws.Range("S21").Value = 0.09855432111742879
$ws.Range("T21").Value = 0.09894431986586213
# Row 22
$ws.Range("G22").Value = 7.135609000000001
$ws.Range("H22").Value = 21.406827
$ws.Range("I22").Value = 0.02688695586661
$ws.Range("J22").Value = 0.0269870280248879
$ws.Range("M22").Value = 181.556244
$ws.Range("N22").Value = 544.668732
$ws.Range("O22").Value = 0.393453292404907
$ws.Range("P22").Value = 0.3935455037432071
$ws.Range("Q22").Value = 1295.514368692596
$ws.Range("R22").Value = 11659.62931823336
$ws.Range("S22").Value = 0.01057876130846313
$ws.Range("T22").Value = 0.01062062353858656
# Row 23
$ws.Range("G23").Value = 7.135609000000001
$ws.Range("H23").Value = 21.406827
$ws.Range("I23").Value = 0.02688695586661
$ws.Range("J23").Value = 0.0269870280248879
$ws.Range("M23").Value = 0.324361
$ws.Range("N23").Value = 0.648722
$ws.Range("O23").Value = 0.0007029276469155644
$ws.Range("P23").Value = 0.0004687282586276696
$ws.Range("Q23").Value = 2.314513270849
$ws.Range("R23").Value = 13.887079625094
$ws.Range("S23").Value = 0.00001889958462003879
$ws.Range("T23").Value = 0.00001264958265164182
# Row 24
$ws.Range("G24").Value = 7.135609000000001
$ws.Range("H24").Value = 21.406827
$ws.Range("I24").Value = 0.02688695586661
$ws.Range("J24").Value = 0.0269870280248879
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 155.929759
$ws.Range("N24").Value = 467.789277
$ws.Range("O24").Value = 0.3379177477501335
$ws.Range("P24").Value = 0.3379969435488647
$ws.Range("Q24").Value = 1112.653791688231
$ws.Range("R24").Value = 10013.88412519408
$ws.Range("S24").Value = 0.009085579570302088
$ws.Range("T24").Value = 0.009121532987879665
# Row 25
$ws.Range("G25").Value = 7.135609000000001
$ws.Range("H25").Value = 21.406827
$ws.Range("I25").Value = 0.02688695586661
$ws.Range("J25").Value = 0.0269870280248879
$ws.Range("M25").Value = 123.632576
$ws.Range("N25").Value = 370.897728
$ws.Range("O25").Value = 0.2679260321980438
$ws.Range("P25").Value = 0.2679888244493004
$ws.Range("Q25").Value = 882.1937219987843
$ws.Range("R25").Value = 7939.743497989058
$ws.Range("S25").Value = 0.007203715403224732
$ws.Range("T25").Value = 0.007232221915770033
